$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Step1 - Input Data")

# B4 holds the "Personal Required Rate of Return" as literal text (e.g. "8.69%"),
# not a numeric percentage. Setting .Value directly would make Excel's COM layer
# auto-convert the string into a numeric percentage, losing the original text
# cell type. Force the cell to text first, assign the new text, then restore
# the original number format so the cell keeps its original style/appearance.
$cell = $ws.Range("B4")
$originalFormat = $cell.NumberFormat
$cell.NumberFormat = "@"
$cell.Value = "8.72%"
$cell.NumberFormat = $originalFormat

$excel.CalculateFullRebuild()
